# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets to match the refreshed data export.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F4").Value = 1497
    $ws.Range("F6").Value = 37
    $ws.Range("F7").Value = 120
    $ws.Range("F9").Value = 295
}
